$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should share the same bold/border
# header style already used by B1:H1. Copy that formatting from H1, then set
# the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill I and J columns for each data row (rows 2 through 26).
# I column is a constant 1, J column mirrors the H column value.
for ($row = 2; $row -le 26; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value2 = 1
    $ws.Cells.Item($row, 10).Value2 = $hValue
}
